# Rule 5's Action cell (D11) was entered with a typo: the string was
# stored as   bad"   instead of the correctly quoted   "bad"   used by the
# equivalent rule in D8. Re-enter it correctly.
#
# A leading apostrophe forces a literal text entry (so Excel/the engine
# doesn't try to treat/evaluate the leading double quote specially) and
# keeps the cell's existing "quote prefix" text formatting intact, just
# like retyping the value by hand would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "'" + '"bad"'

# After correcting D11 and pressing Enter, the active cell moves down to D12.
$ws.Range("D12").Select() | Out-Null
